$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 18.8157901763916
$ws.Range("C2").Value = 5.344827651977539
$ws.Range("D2").Value = 12.418232917785645
$ws.Range("E2").Value = 57.85714340209961
